$d = $word.ActiveDocument

# ------------------------------------------------------------------
# This schedule document had its "Oct 23 - Nov 18" block of entries
# copied earlier in the list (replacing the stale "Oct 14 - Oct 21"
# block), a couple of leftover duplicate lines removed, and the
# final week's sub-bullets reworked from a quiz-topic outline into a
# simple Quiz / Review quiz / Review Quiz set of bullets.
#
# We apply the operations working from the bottom of the document
# upward so that paragraph indices for not-yet-processed items never
# shift underneath us.
# ------------------------------------------------------------------

# 1) After "Monday 11/18/2024" (para 59), add a new sub-bullet
#    "Review Quiz" at the second outline level (ilvl = 1).
$pMonday1118 = $d.Paragraphs.Item(59)
if ($pMonday1118.Range.Text.TrimEnd([char]13) -ne "Monday 11/18/2024") {
    throw "Unexpected paragraph 59: " + $pMonday1118.Range.Text
}
$pMonday1118.Range.InsertParagraphAfter()
$pReviewQuiz = $d.Paragraphs.Item(60)
$pReviewQuiz.Range.Text = "Review Quiz"
$pReviewQuiz.Range.ListFormat.ListLevelNumber = 2

# 2) Turn the old quiz-review sub-bullets (paras 56-58: "AMA Review
#    Session" / "Homework 7 on Chapter 8 due by 6pm" / "Survey")
#    into just two bullets: "Quiz" and "Review quiz".
$p56 = $d.Paragraphs.Item(56)
if ($p56.Range.Text.TrimEnd([char]13) -ne "AMA Review Session") {
    throw "Unexpected paragraph 56: " + $p56.Range.Text
}
$p56.Range.Text = "Quiz"

$p57 = $d.Paragraphs.Item(57)
if ($p57.Range.Text.TrimEnd([char]13) -ne "Homework 7 on Chapter 8 due by 6pm") {
    throw "Unexpected paragraph 57: " + $p57.Range.Text
}
$p57.Range.Text = "Review quiz"

$p58 = $d.Paragraphs.Item(58)
if ($p58.Range.Text.TrimEnd([char]13) -ne "Survey") {
    throw "Unexpected paragraph 58: " + $p58.Range.Text
}
$p58.Range.Delete()

# 3) Remove the stray duplicate "Assign Homework 6 covering Chapter
#    7??" bullet that appeared right before "Monday 11/04/2024".
$p40 = $d.Paragraphs.Item(40)
if ($p40.Range.Text.TrimEnd([char]13) -ne "Assign Homework 6 covering Chapter 7??") {
    throw "Unexpected paragraph 40: " + $p40.Range.Text
}
$p40.Range.Delete()

# 4) Delete the now-superseded "Monday 10/14/2024" ... "Assign
#    Homework 4 covering Chapter 5" block (original paragraphs 2-29)
#    since that content has effectively moved to the start of the
#    list (the "Wednesday 10/23/2024" ... block that already follows
#    it keeps its own, already-correct, copies of this material).
$startPara = $d.Paragraphs.Item(2)
$endPara = $d.Paragraphs.Item(29)
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()

Write-Host ("Final paragraph count: " + $d.Paragraphs.Count)
